$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 169, shifting rows 169:256 down to 170:257
$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value = 5
$ws.Cells.Item(169, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(169, 3).Value = "Maule"
$ws.Cells.Item(169, 4).Value = 44572
$ws.Cells.Item(169, 5).Value = 7
$ws.Cells.Item(169, 6).Value = 100112023
$ws.Cells.Item(169, 7).Value = "Brócoli"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 3000
$ws.Cells.Item(169, 11).Value = 500
$ws.Cells.Item(169, 12).Value = 500
$ws.Cells.Item(169, 13).Value = 500
$ws.Cells.Item(169, 14).Value = "`$/unidad"
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 500
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
